{"js": "// Office.js (Word JavaScript API)\n// Append, at the very end of the document body (after the last table and\n// its trailing empty paragraph, right before the final section break), a\n// blank paragraph, a \"TABLA DE AMORTIZACI\u00d3N\" heading paragraph and a\n// \"{{TABLA_AMORTIZACION}}\" placeholder paragraph.\n\nconst body = context.document.body;\n\n// Insert in order at the end of the body so the final body order becomes:\n//   ... <last existing empty paragraph> <blank> <TABLA DE AMORTIZACI\u00d3N> <{{TABLA_AMORTIZACION}}>\nbody.insertParagraph(\"\", Word.InsertLocation.end);\nbody.insertParagraph(\"TABLA DE AMORTIZACI\u00d3N\", Word.InsertLocation.end);\nbody.insertParagraph(\"{{TABLA_AMORTIZACION}}\", Word.InsertLocation.end);\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style)\n# Append, at the very end of the document (after the last table and its\n# trailing empty paragraph, right before the final section break), a blank\n# paragraph, a \"TABLA DE AMORTIZACI\u00d3N\" heading paragraph and a\n# \"{{TABLA_AMORTIZACION}}\" placeholder paragraph.\n\n$d = $word.ActiveDocument\n\n# 1) Append one blank paragraph at the very end of the document.\n$r = $d.Content\n$r.Collapse(0)\n$r.InsertParagraphAfter()\n\n# 2) Append the \"TABLA DE AMORTIZACI\u00d3N\" heading and the\n#    \"{{TABLA_AMORTIZACION}}\" placeholder as two further paragraphs (a\n#    leading paragraph mark keeps the blank paragraph from step 1 empty).\n$r = $d.Content\n$r.Collapse(0)\n$r.InsertAfter(\"`rTABLA DE AMORTIZACI\u00d3N`r{{TABLA_AMORTIZACION}}\")\n"}
